# fix converter for answers
# also shorten tab name for DGA formulaire

$wb = $excel.ActiveWorkbook

# --- Rename the worksheet tabs -------------------------------------------
# "Formulaire-SDI-SecNum-2216" -> "Formulaire"
# "answers"                    -> "réponses"
$wb.Worksheets.Item("Formulaire-SDI-SecNum-2216").Name = "Formulaire"
$wb.Worksheets.Item("answers").Name = "réponses"

$wsLib  = $wb.Worksheets.Item("library_content")
$wsForm = $wb.Worksheets.Item("Formulaire")
$wsRep  = $wb.Worksheets.Item("réponses")

# --- Update the library_content rows that record the sheet/tab names ----
# (row 14 -> Formulaire tab ref, row 15 -> réponses tab ref)
$wsLib.Cells.Item(14, 2).Value = "Formulaire"
$wsLib.Cells.Item(15, 2).Value = "réponses"

# --- Restore each sheet's own selection, then finish on library_content -
$wsForm.Activate()
$wsForm.Range("D39").Select()

$wsRep.Activate()
$wsRep.Range("C3").Select()

$wsLib.Activate()
$wsLib.Range("B16").Select()
